$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two observation records that occupied rows 8 and 9 traded places:
# everything that is specific to the record (id, taxon id/name, K/L/M/N
# activity fields, coordinates, ...) swaps between the rows, while the
# columns that already held identical data in both rows (C, D, I, P, S,
# T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY, ...) stay put.

# Columns that hold a value in both row 8 and row 9 today - simple swap.
$swapCols = @("A","B","E","F","G","H","Q","R")

foreach ($col in $swapCols) {
    $top = $ws.Range("$col" + "8")
    $bot = $ws.Range("$col" + "9")
    $topVal = $top.Value()
    $botVal = $bot.Value()
    $top.Value = $botVal
    $bot.Value = $topVal
}

# Columns K, L, M, N currently only carry data on row 9 (row 8 is blank).
# After the swap, row 8 gets that data (K8/L8/N8 are empty-text cells,
# same as K9/L9/N9 are today) and row 9 becomes blank instead.
# A leading "'" forces literal/empty text (same convention Excel itself
# uses for a cell typed as just an apostrophe) instead of collapsing an
# assigned "" down to a truly blank cell; the style is reset straight
# back to Normal afterwards so the quote-prefix formatting it implies
# doesn't stick around.
foreach ($col in @("K","L","N")) {
    $cell = $ws.Range("$col" + "8")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

$ws.Range("M8").Value = $ws.Range("M9").Value()

$ws.Range("K9").Value = $null
$ws.Range("L9").Value = $null
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = $null
